$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column F (bold, matching the other header cells)
$ws.Range("F1").Value = "Bán trưc tiếp"
$ws.Range("F1").Font.Bold = $true

# Set column F width to match diff (bestFit width ~12.2)
$ws.Columns.Item(6).ColumnWidth = 11.428571428571429

# Fill in new data column values
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1

# Update selection to reflect new active cell
$ws.Range("F4").Select()
